# Adds the "2022-Q1" sheet (new holdings detail, inserted right before "总计")
# and updates the "总计" (summary) sheet with a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet right before "总计"
# ---------------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$ws = $wb.Worksheets.Add($zongji)
$ws.Name = "2022-Q1"

# Header row (row 1) -------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$cols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Range($cols[$i] + "1")
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data rows (rows 2-14) ------------------------------------------------------
$rows = @(
    @("003567", "华夏行业景气混合",           "112.21", "91.63", "2.51", "2.8165", 10),
    @("001245", "工银瑞信生态环境行业股票",    "57.90",  "90.05", "3.38", "1.9570", 8),
    @("009147", "建信新能源行业股票",          "54.78",  "87.41", "3.41", "1.8680", 9),
    @("013175", "海富通碳中和混合A",           "14.27",  "93.69", "3.30", "0.4709", 9),
    @("530001", "建信恒久价值混合",            "11.95",  "91.57", "2.94", "0.3513", 9),
    @("013176", "海富通碳中和混合C",           "7.99",   "93.69", "3.30", "0.2637", 9),
    @("000592", "建信改革红利股票",            "9.09",   "91.05", "2.56", "0.2327", 10),
    @("008177", "建信高股息主题股票",          "5.58",   "93.26", "2.90", "0.1618", 7),
    @("001858", "建信鑫利灵活配置混合",        "5.41",   "85.72", "2.61", "0.1412", 9),
    @("290014", "泰信现代服务业混合",          "0.73",   "81.14", "6.44", "0.0470", 8),
    @("005009", "申万菱信行业轮动股票",        "0.61",   "90.25", "4.33", "0.0264", 8),
    @("007965", "民生加银品质消费股票A",       "0.25",   "88.39", "5.48", "0.0137", 4),
    @("007966", "民生加银品质消费股票C",       "0.13",   "88.39", "5.48", "0.0071", 4)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $data = $rows[$r]

    $aCell = $ws.Range("A" + $rowNum)
    $aCell.Value = $r
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    # B: fund code, C: fund name -> plain text
    $ws.Range("B" + $rowNum).NumberFormat = "@"
    $ws.Range("B" + $rowNum).Value = $data[0]
    $ws.Range("C" + $rowNum).Value = $data[1]

    # D-G: numeric-looking values stored as plain text (matches source data)
    $ws.Range("D" + $rowNum).NumberFormat = "@"
    $ws.Range("D" + $rowNum).Value = $data[2]
    $ws.Range("E" + $rowNum).NumberFormat = "@"
    $ws.Range("E" + $rowNum).Value = $data[3]
    $ws.Range("F" + $rowNum).NumberFormat = "@"
    $ws.Range("F" + $rowNum).Value = $data[4]
    $ws.Range("G" + $rowNum).NumberFormat = "@"
    $ws.Range("G" + $rowNum).Value = $data[5]

    # H: rank -> real number
    $ws.Range("H" + $rowNum).Value = $data[6]
}

$ws.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2. Update "总计" - insert a new leading row for the 2022-Q1 aggregate
# ---------------------------------------------------------------------------
# NOTE: re-fetch the "总计" worksheet by name. The variable captured before
# Worksheets.Add/Name (re)binds to whatever sheet is active afterwards, so a
# stale reference here would silently write into the wrong sheet.
$zongji = $wb.Worksheets.Item("总计")
$zongji.Rows.Item(2).Insert()

$a2 = $zongji.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

$zongji.Range("B2:D2").ClearFormats()
$zongji.Range("B2").Value = "2022-Q1"
$zongji.Range("C2").Value = 13
$zongji.Range("D2").Value = 8.36

# The existing rows (now 3-7) keep their old index label in column A from
# before the insert (0,1,2,3,4) - renumber them sequentially (1,2,3,4,5) so
# the whole A2:A7 run reads 0,1,2,3,4,5 top-to-bottom.
for ($r = 3; $r -le 7; $r++) {
    $zongji.Range("A" + $r).Value = $r - 2
}

$zongji.Range("A1").Select()
